# Generate Report for Handoff
#
# Refresh the localization-status report with a fresh handoff timestamp for
# the "985a6e46-20a0-4ec2-b4e8-f00f7999be14.md" file (row 7 on each sheet):
#   - Overview!D7            (Latest Handoff Date)
#   - zh-cn!E7                (Latest Handoff Datetime)
#   - de-de!E7                (Latest Handoff Datetime)

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("D7").Value = "2016-03-21 16:41:08"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("E7").Value = "2016-03-21 16:41:01"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("E7").Value = "2016-03-21 16:41:08"
